$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Data for the nine new rows (5..13): they cycle through the same three
#    "item" templates already used by rows 2/3/4, only the per-row Number
#    and ImageUrl change.
# ---------------------------------------------------------------------------
$rows = @(
    @{ r = 5;  n = 4;  item = 1 },
    @{ r = 6;  n = 5;  item = 2 },
    @{ r = 7;  n = 6;  item = 3 },
    @{ r = 8;  n = 7;  item = 1 },
    @{ r = 9;  n = 8;  item = 2 },
    @{ r = 10; n = 9;  item = 3 },
    @{ r = 11; n = 10; item = 1 },
    @{ r = 12; n = 11; item = 2 },
    @{ r = 13; n = 12; item = 3 }
)

$titlesEn = @{ 1 = "Item 1"; 2 = "Item 2"; 3 = "Item 3" }
$descEn = @{
    1 = "This is the description of item 1."
    2 = "This is the description of item 2."
    3 = "This is the description of item 3."
}
$titlesAr = @{ 1 = "عنصر 1"; 2 = "عنصر 2"; 3 = "عنصر 3" }
$descAr = @{
    1 = "هذا هو وصف العنصر 1."
    2 = "هذا هو وصف العنصر 2."
    3 = "هذا هو وصف العنصر 3."
}
$price = @{ 1 = 10; 2 = 15; 3 = 20 }

# ---------------------------------------------------------------------------
# 2. Write the values first (so the sheet grows to A1:G13) ...
# ---------------------------------------------------------------------------
foreach ($row in $rows) {
    $r = $row.r
    $it = $row.item
    $ws.Cells.Item($r, 1).Value = $row.n
    $ws.Cells.Item($r, 2).Value = $titlesEn[$it]
    $ws.Cells.Item($r, 3).Value = $descEn[$it]
    $ws.Cells.Item($r, 4).Value = $price[$it]
    $ws.Cells.Item($r, 5).Value = "https://raw.githubusercontent.com/dofreelancer19/showitems/main/images/$($row.n).jpeg"
    $ws.Cells.Item($r, 6).Value = $titlesAr[$it]
    $ws.Cells.Item($r, 7).Value = $descAr[$it]
}

# ---------------------------------------------------------------------------
# 3. ... then hyperlink column E for each new row (before the formatting
#    copy below, since Hyperlinks.Add mints its own style - applying the
#    format copy afterwards restores the exact shared style the pasted
#    rows use, just like the original file). Added row-by-row in the same
#    order the author's own edit shows in the saved file (two rows from
#    each pasted 3-row block, then the three rows filled in afterwards).
# ---------------------------------------------------------------------------
$byRow = @{}
foreach ($row in $rows) { $byRow[$row.r] = $row }
$hyperlinkOrder = @(6, 7, 9, 10, 12, 13, 5, 8, 11)
foreach ($r in $hyperlinkOrder) {
    $row = $byRow[$r]
    $cell = $ws.Cells.Item($r, 5)
    $ws.Hyperlinks.Add($cell, "https://raw.githubusercontent.com/dofreelancer19/showitems/main/images/$($row.n).jpeg")
}

# ---------------------------------------------------------------------------
# 4. Copy the formatting (styles + row height) of the existing "wrapped"
#    data rows (3:4 -> style s=3/4/5, 41.4pt row height, hyperlink style on
#    column E) down across the new rows 5:13. PasteSpecial xlPasteFormats
#    (-4122) tiles the 2-row source block across the 9-row destination,
#    reproducing the same alternation already used by rows 2/3/4.
# ---------------------------------------------------------------------------
$ws.Range("A3:G4").Copy()
$ws.Range("A5:G13").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row height: the source rows used an auto-computed 41.4pt height (wrapped
# text); reproduce the same value explicitly for all new rows.
$ws.Range("A5:G13").RowHeight = 41.4

# ---------------------------------------------------------------------------
# 5. Selection left wherever the author's last click landed after pasting.
# ---------------------------------------------------------------------------
$ws.Range("I7").Select()
